$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the contribution percentages (D4:D7) for the work-assignment table
$ws.Range("D4").Value = 0.4
$ws.Range("D5").Value = 0.2
$ws.Range("D6").Value = 0.2
$ws.Range("D7").Value = 0.2

# Match the saved selection state from the edit
$ws.Range("E8").Select()
